$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values for each data row (G stays 1, H stays blank)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (row 2 .. row 25), in column order above
$data = @{
    2 = @(1.02, 1.04392369179643, 1.043364074513166, 1.050832353931925, 1.059253070185213, 1.042512664174057, 1.048991534437422, 1.046138149082773, 1.053585536458262, 1.061983079034576, 1.020131371447747)
    3 = @(1.02, 1.045133271871064, 1.044020259616436, 1.051945613268927, 1.06053905461577, 1.042818860456435, 1.049846871600676, 1.046605759135271, 1.054510546737093, 1.063082062818689, 1.020421704973886)
    4 = @(1.02, 1.045915515579285, 1.044444316176806, 1.052665948803083, 1.061371442928274, 1.043015292385641, 1.050399382343965, 1.046907111971836, 1.055108484987343, 1.063792894855228, 1.020609090317929)
    5 = @(1.02, 1.046244268739252, 1.044622459960051, 1.052968775026471, 1.061721445578018, 1.043097466150268, 1.050631431758793, 1.047033507948724, 1.055359714626285, 1.064091662267457, 1.020687752561006)
    6 = @(1.02, 1.046299461917217, 1.044652363469797, 1.053019620756293, 1.061780216522937, 1.043111239671117, 1.050670380653661, 1.047054713216884, 1.055401888804906, 1.064141822784857, 1.020700953589552)
    7 = @(1.02, 1.045919908793571, 1.044446697052676, 1.0526699951913, 1.061376119416775, 1.04301639199146, 1.050402483886502, 1.046908802031633, 1.055111842492721, 1.063796887258831, 1.020610141856395)
    8 = @(1.02, 1.044332566517252, 1.043585946640835, 1.051208589647129, 1.05968761978449, 1.042616496235983, 1.049280797000884, 1.04629643303871, 1.053898274117211, 1.062354545512511, 1.020229590576506)
    9 = @(1.02, 1.041532034184943, 1.042065088661523, 1.048633207495822, 1.056714254407486, 1.041898819251203, 1.047296912838042, 1.045208007454185, 1.051755117469668, 1.059810707972165, 1.019555320693144)
    10 = @(1.02, 1.039662566536752, 1.041048451228609, 1.046916061135404, 1.054733226213204, 1.041411607306655, 1.045969315643255, 1.044476108725169, 1.050323105981294, 1.058113201806891, 1.019103304961541)
    11 = @(1.02, 1.038852454672914, 1.040607593155805, 1.046172443318024, 1.053875679804079, 1.041198557590277, 1.045393244960928, 1.044157699222684, 1.04970224110214, 1.057377756109588, 1.018906978491213)
    12 = @(1.02, 1.038551447083512, 1.040443741906396, 1.045896216294295, 1.053557184249825, 1.041119107983828, 1.045179082812338, 1.044039203475043, 1.049471503107073, 1.057104514430049, 1.018833963372751)
    13 = @(1.02, 1.038616018621598, 1.040478892936867, 1.045955468604645, 1.053625501066931, 1.04113616438511, 1.045225029667211, 1.044064631365193, 1.049521002700508, 1.057163128646855, 1.018849629476815)
    14 = @(1.02, 1.038827575265838, 1.040594051135955, 1.046149610608047, 1.053849352155145, 1.041191996657753, 1.045375546011263, 1.044147908902789, 1.049683170708912, 1.057355171193388, 1.018900944892319)
    15 = @(1.02, 1.038957909494743, 1.040664991068235, 1.046269225957079, 1.053987278810492, 1.041226355215667, 1.04546825968567, 1.044199189221224, 1.049783071666259, 1.057473486346058, 1.018932549976923)
    16 = @(1.02, 1.039716317689242, 1.041077695916018, 1.046965410707186, 1.054790143839594, 1.041425702778963, 1.046007521869096, 1.044497209065678, 1.050364293848806, 1.058162002002658, 1.019116321806235)
    17 = @(1.02, 1.040191879123645, 1.041336401681378, 1.047402085670814, 1.0552938255481, 1.041550190085094, 1.0463454610094, 1.044683749419313, 1.050728665217579, 1.058593777238362, 1.019231435862619)
    18 = @(1.02, 1.040469206430176, 1.041487237896755, 1.047656782997913, 1.055587639038045, 1.041622600402846, 1.046542458336429, 1.044792411249284, 1.050941120340433, 1.05884558459895, 1.019298522020684)
    19 = @(1.02, 1.040563757844924, 1.041538658539963, 1.047743626940445, 1.055687826070278, 1.041647256364906, 1.0466096095644, 1.044829437700871, 1.051013549119085, 1.058931437710452, 1.019321386856461)
    20 = @(1.02, 1.040140862072076, 1.041308651461977, 1.047355235400691, 1.055239782755544, 1.041536854571217, 1.0463092154381, 1.044663750300425, 1.050689579547932, 1.058547455978251, 1.019219091208887)
    21 = @(1.02, 1.038765279779038, 1.040560142589804, 1.046092441035932, 1.053783432632715, 1.041175564101188, 1.045331227788314, 1.044123391935282, 1.049635419630599, 1.057298621245118, 1.018885836289815)
    22 = @(1.02, 1.037899839921419, 1.04008896468802, 1.045298388655918, 1.052867970882281, 1.040946592683149, 1.044715263599075, 1.043782349532865, 1.048971925715888, 1.05651305560017, 1.018675780542675)
    23 = @(1.02, 1.038358679675328, 1.040338798047688, 1.045719339399837, 1.053353256029522, 1.041068146822435, 1.045041899348887, 1.04396326558342, 1.049323723505961, 1.056929534954568, 1.018787185008531)
    24 = @(1.02, 1.040163914665907, 1.041321190784898, 1.047376405045241, 1.05526420229167, 1.041542880937104, 1.046325593612692, 1.044672787486477, 1.050707240916126, 1.058568386684911, 1.019224669405539)
    25 = @(1.02, 1.042256459892462, 1.042458750491042, 1.049299038442541, 1.057482717557187, 1.04208589904563, 1.047810671172111, 1.045490499213969, 1.052309739767059, 1.060468627888491, 1.019730075106251)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
